$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 edits ---
# E6: clear the existing value "4", leaving the cell blank but still present
$ws.Range("E6").Value = ""
$ws.Range("E6").Style = "Normal"
# G6: change "6" to "no pero si"
$ws.Range("G6").Value = "no pero si"

# --- Row 8 edits ---
# C8, E8, F8, H8 were already blank placeholder cells; clear them so they
# are dropped entirely (D8 and G8 keep their existing values untouched)
$ws.Range("C8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("H8").Value = ""

# --- New row 9 ---
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "67555"
$ws.Range("A9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3"
$ws.Range("E9").Style = "Normal"

$ws.Range("G9").Value = "a"

# --- New row 10 ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "12345"
$ws.Range("A10").Style = "Normal"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "1"
$ws.Range("B10").Style = "Normal"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "2"
$ws.Range("C10").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1"
$ws.Range("E10").Style = "Normal"

$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "5"
$ws.Range("F10").Style = "Normal"

$ws.Range("G10").Value = "COSINES"

$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "1"
$ws.Range("H10").Style = "Normal"
